$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.513.01'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.913.53'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = '''0.722'
$ws.Range('E5').Value = '  +10.12%  '
$ws.Range('D6').Value = '''248.05'
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('D8').Value = '''40.65'
$ws.Range('E8').Value = '  -3.98%  '
$ws.Range('E9').Value = '  +2.07%  '
$ws.Range('D10').Value = '''53.10'
$ws.Range('E10').Value = '  +7.88%  '
$ws.Range('D11').Value = '''0.0735'
$ws.Range('E11').Value = '  +2.37%  '
$ws.Range('D12').Value = '''0.0988'
$ws.Range('E12').Value = '  -1.69%  '
$ws.Range('D13').Value = '2.189.97'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').Value = '''12.62'
$ws.Range('E14').Value = '  +2.65%  '
$ws.Range('D15').Value = '''0.718'
$ws.Range('E15').Value = '  +2.20%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.916.03'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '''4.91'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '35.514.31'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = '''73.13'
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '''13.19'
$ws.Range('E21').Value = '  +3.72%  '
$ws.Range('D22').Value = '''242.22'
$ws.Range('E22').Value = '  -1.46%  '
$ws.Range('D23').Value = '''5.06'
$ws.Range('E23').Value = '  +4.36%  '
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').Value = '''2.36'
$ws.Range('E26').Value = '  +7.33%  '
$ws.Range('D27').Value = '''167.69'
$ws.Range('E27').Value = '  -2.31%  '
$ws.Range('D28').Value = '''8.67'
$ws.Range('E28').Value = '  +2.09%  '
$ws.Range('D29').Value = '''18.84'
$ws.Range('E29').Value = '  +0.67%  '
$ws.Range('D30').Value = '''0.134'
$ws.Range('E30').Value = '  +3.94%  '
$ws.Range('D31').Value = '4.142.86'
$ws.Range('E31').Value = '  +19.87%  '
$ws.Range('D32').Value = '''4.37'
$ws.Range('E32').Value = '  +5.09%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '''1.97'
$ws.Range('E33').Value = '  +13.46%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '''0.0578'
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('E35').Value = '  +0.73%  '
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').Value = '''0.917'
$ws.Range('E37').Value = '  -5.24%  '
$ws.Range('E38').Value = '  +9.83%  '
$ws.Range('D39').Value = '''2.05'
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').Value = '''17.97'
$ws.Range('E40').Value = '  +13.77%  '
$ws.Range('D41').Value = '''98.91'
$ws.Range('E41').Value = '  +6.94%  '
$ws.Range('D42').Value = '''1.14'
$ws.Range('E42').Value = '  +1.74%  '
$ws.Range('D43').Value = '''0.0210'
$ws.Range('E43').Value = '  -0.92%  '
$ws.Range('D44').Value = '''0.0655'
$ws.Range('E44').Value = '  +3.93%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '''2.50'
$ws.Range('E45').Value = '  +2.27%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.346.86'
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('D47').Value = '''2.43'
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('E48').Value = '  -0.72%  '
$ws.Range('D49').Value = '''45.42'
$ws.Range('E49').Value = '  -3.49%  '
$ws.Range('E50').Value = '  -0.33%  '
$ws.Range('D51').Value = '''12.06'
$ws.Range('E51').Value = '  -6.60%  '
